# 🔄 MAJ automatique BRVM via GitHub Actions
# Refresh of the "Recommandations" and "Top_YTD" sheets with the latest
# BRVM market figures (variation totals, last variation, YTD progress)
# and the resulting re-ranking of tickers within each sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Recommandations sheet ---
$ws1.Range("D2").Value = 3254.93
$ws1.Range("E2").Value = 104.49
$ws1.Range("D3").Value = 2930
$ws1.Range("E3").Value = 800
$ws1.Range("D4").Value = 2612.38
$ws1.Range("E4").Value = 643.52
$ws1.Range("D5").Value = 2605
$ws1.Range("E5").Value = 650
$ws1.Range("A6").Value = "SETAO CI"
$ws1.Range("D6").Value = 2320
$ws1.Range("E6").Value = 580
$ws1.Range("A7").Value = "NEI-CEDA CI"
$ws1.Range("D7").Value = 2295
$ws1.Range("E7").Value = 565
$ws1.Range("D8").Value = 2260
$ws1.Range("E8").Value = 575
$ws1.Range("D9").Value = 2060
$ws1.Range("E9").Value = 505
$ws1.Range("D10").Value = 1433.2
$ws1.Range("E10").Value = 359.89
$ws1.Range("D11").Value = 1386.63
$ws1.Range("E11").Value = 352.47
$ws1.Range("D12").Value = 1247.55
$ws1.Range("E12").Value = 316.63
$ws1.Range("D13").Value = 993.27
$ws1.Range("E13").Value = 255.28
$ws1.Range("D14").Value = 822.64
$ws1.Range("E14").Value = 211.08
$ws1.Range("D15").Value = 742.51
$ws1.Range("E15").Value = 187.08
$ws1.Range("D16").Value = 545.17
$ws1.Range("E16").Value = 138.83
$ws1.Range("D17").Value = 518.51
$ws1.Range("E17").Value = 128.86
$ws1.Range("D18").Value = 489.75
$ws1.Range("E18").Value = 121.53
$ws1.Range("D19").Value = 481.33
$ws1.Range("E19").Value = 119.44
$ws1.Range("D20").Value = 425.42
$ws1.Range("E20").Value = 105
$ws1.Range("D21").Value = 418.69
$ws1.Range("E21").Value = 105.9
$ws1.Range("D22").Value = 368.87
$ws1.Range("E22").Value = 92.05
$ws1.Range("E23").Value = 7.5
$ws1.Range("A24").Value = "SAFCA CI (SAFC)"
$ws1.Range("D24").Value = 14.57
$ws1.Range("E24").Value = 7.38
$ws1.Range("A25").Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws1.Range("D25").Value = 6.33
$ws1.Range("E25").Value = 6.33
$ws1.Range("A26").Value = "SAPH CI (SPHC)"
$ws1.Range("B26").Value = 1
$ws1.Range("D26").Value = 4.51
$ws1.Range("E26").Value = 4.51
$ws1.Range("A27").Value = "ORAGROUP TOGO (ORGT)"
$ws1.Range("C27").Value = 1
$ws1.Range("D27").Value = 4.09
$ws1.Range("E27").Value = -1.74
$ws1.Range("G27").Value = "👀 À surveiller"
$ws1.Range("A28").Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws1.Range("D28").Value = 3.98
$ws1.Range("E28").Value = -3.42
$ws1.Range("A29").Value = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$ws1.Range("B29").Value = 2
$ws1.Range("C29").Value = 1
$ws1.Range("D29").Value = 2.88
$ws1.Range("E29").Value = -3.03
$ws1.Range("G29").Value = "👀 À surveiller"
$ws1.Range("A32").Value = "BERNABE CI (BNBC)"
$ws1.Range("B32").Value = 2
$ws1.Range("D32").Value = 1.18
$ws1.Range("E32").Value = 3.59
$ws1.Range("A33").Value = "CIE CI (CIEC)"
$ws1.Range("D33").Value = 0.51
$ws1.Range("E33").Value = -5.06
$ws1.Range("A34").Value = "BANK OF AFRICA BN (BOAB)"
$ws1.Range("D34").Value = 0.42
$ws1.Range("E34").Value = 2.86
$ws1.Range("A36").Value = "FILTISAC CI (FTSC)"
$ws1.Range("D36").Value = -0.61
$ws1.Range("E36").Value = 2.89
$ws1.Range("A41").Value = "SITAB CI (STBC)"
$ws1.Range("D41").Value = -2.5
$ws1.Range("E41").Value = -2.5
$ws1.Range("A42").Value = "SICOR CI (SICC)"
$ws1.Range("D42").Value = -2.73
$ws1.Range("E42").Value = -2.73
$ws1.Range("A43").Value = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$ws1.Range("D43").Value = -3.46
$ws1.Range("E43").Value = -3.46
$ws1.Range("A44").Value = "CFAO MOTORS CI (CFAC)"
$ws1.Range("B44").Value = 0
$ws1.Range("C44").Value = 1
$ws1.Range("D44").Value = -3.65
$ws1.Range("E44").Value = -3.65
$ws1.Range("G44").Value = "➖ Neutre"
$ws1.Range("A45").Value = "SONATEL SN (SNTS)"
$ws1.Range("D45").Value = -3.81
$ws1.Range("E45").Value = -3.81
$ws1.Range("A46").Value = "SUCRIVOIRE (SCRC)"
$ws1.Range("C46").Value = 1
$ws1.Range("D46").Value = -3.85
$ws1.Range("E46").Value = -3.85
$ws1.Range("A47").Value = "BANK OF AFRICA NG (BOAN)"
$ws1.Range("D47").Value = -7.05
$ws1.Range("E47").Value = -2

# --- Top_YTD sheet ---
$ws2.Range("B2").Value = 7644798.21
$ws2.Range("B3").Value = 477532.03
$ws2.Range("B4").Value = 321526.5
$ws2.Range("B5").Value = 318246.88
$ws2.Range("A6").Value = "SETAO CI"
$ws2.Range("B6").Value = 213702.2
$ws2.Range("A7").Value = "NEI-CEDA CI"
$ws2.Range("B7").Value = 205694.23
$ws2.Range("B8").Value = 195319.25
$ws2.Range("B9").Value = 142878.52
$ws2.Range("B10").Value = 44015.82
$ws2.Range("B11").Value = 39696.24
